$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the list of names in column A (rows 2-99) alphabetically (A to Z),
# the same way Excel's Data > Sort command would, which also records a
# <sortState> on the worksheet.
$sortRange = $ws.Range("A2:A99")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A99"))
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 2
$ws.Sort.Apply()

# Move the selection back to D5 (and, implicitly, scroll the view back to
# the top so the saved view no longer shows topLeftCell="A82").
$ws.Range("D5").Select()
